$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''25.866.72'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '''1.629.14'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '''214.50'
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').Value = '''0.502'
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').Value = '''0.255'
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').Value = '''0.0631'
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').Value = '''19.64'
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('D11').Value = '''0.0787'
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('D12').Value = '''1.852.29'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').Value = '''4.24'
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('D14').Value = '''1.615.50'
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('D15').Value = '''0.544'
$ws.Range('E15').Value = '  -1.90%  '
$ws.Range('D16').Value = '''0.0₃0758'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').Value = '''62.73'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('D18').Value = '''25.865.49'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('E19').Value = '  -0.10%  '
$ws.Range('D20').Value = '''193.00'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').Value = '''4.38'
$ws.Range('E21').Value = '  -1.04%  '
$ws.Range('D22').Value = '''9.96'
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('D23').Value = '''6.24'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('E24').Value = '  -1.83%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '''142.67'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('B26').Value = 'BinanceUSD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D26').Value = '''0.998'
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('D27').Value = '''0.125'
$ws.Range('E27').Value = '  +1.69%  '
$ws.Range('D28').Value = '''6.86'
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = '''15.43'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').Value = '''1.24'
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').Value = '''0.0498'
$ws.Range('E31').Value = '  +2.03%  '
$ws.Range('D32').Value = '''3.31'
$ws.Range('E32').Value = '  -0.57%  '
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('E34').Value = '  +0.66%  '
$ws.Range('E35').Value = '  +1.83%  '
$ws.Range('D36').Value = '''0.900'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').Value = '''1.133.98'
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('D38').Value = '''0.547'
$ws.Range('E38').Value = '  +1.04%  '
$ws.Range('E39').Value = '  -2.20%  '
$ws.Range('D40').Value = '''0.0156'
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = '''99.36'
$ws.Range('E42').Value = '  -0.72%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '''5.46'
$ws.Range('E43').Value = '  -0.84%  '
$ws.Range('D44').Value = '''0.799'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').Value = '''1.762.74'
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('D47').Value = '''56.12'
$ws.Range('E47').Value = '  +1.65%  '
$ws.Range('D48').Value = '''0.0528'
$ws.Range('E48').Value = '  +3.00%  '
$ws.Range('D49').Value = '''1.45'
$ws.Range('E49').Value = '  +2.11%  '
$ws.Range('D50').Value = '''0.415'
$ws.Range('E50').Value = '  -0.87%  '
$ws.Range('E51').Value = '  +2.76%  '
